$wb = $excel.ActiveWorkbook

# Sheet: zh-cn -> new handoff datetime for files just handed off
$wsZh = $wb.Worksheets.Item("zh-cn")
$zhHandoffTime = "2016-03-02 10:47:38"
for ($r = 2; $r -le 11; $r++) {
    $status = $wsZh.Cells.Item($r, 2).Text
    if ($status -eq "Ready for handoff" -or $status -eq "Handback transform failed") {
        $wsZh.Cells.Item($r, 4).Value = $zhHandoffTime
    }
}

# Sheet: de-de -> new handoff datetime for files just handed off
$wsDe = $wb.Worksheets.Item("de-de")
$deHandoffTime = "2016-03-02 10:47:50"
for ($r = 2; $r -le 11; $r++) {
    $status = $wsDe.Cells.Item($r, 2).Text
    if ($status -eq "Ready for handoff" -or $status -eq "Handback transform failed") {
        $wsDe.Cells.Item($r, 4).Value = $deHandoffTime
    }
}
